# Add a new "2022" data column (column J) to the table, mirroring the
# existing "2021" column (I) for both formatting and layout, then fill in
# the new figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone column I's formatting (borders/number-format/font) for rows 3-14
# into column J, so the new column matches the look of the rest of the
# table before we overwrite it with the 2022 values.
$ws.Range("I3:I14").Copy()
$ws.Range("J3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New 2022 figures.
$ws.Range("J4").Value = 2022
$ws.Range("J5").Value = 96.4
$ws.Range("J6").Value = 96.4
$ws.Range("J7").Value = 97.9
$ws.Range("J8").Value = 95.3
$ws.Range("J9").Value = 93.8
$ws.Range("J10").Value = 95.5
$ws.Range("J11").Value = 94.4
$ws.Range("J12").Value = 95
$ws.Range("J13").Value = 98.7
$ws.Range("J14").Value = 97.3

# Match the saved selection from the source workbook.
[void]$ws.Range("L10").Select()
